$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Row 10 ("hidden" / "hidden_field" / "Hidden field test") is being removed
# from the middle of the survey sheet and re-added as a new row at the very
# end of the sheet (the xlsform test fixture wants the hidden-field question
# to appear last). Deleting the whole row shifts every following row up by
# one, which matches the rest of the observed diff.
$ws.Rows.Item(10).Delete()

# Re-append the same question (type / name / label) as the new last row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1
$ws.Cells.Item($lastRow, 1).Value = "hidden"
$ws.Cells.Item($lastRow, 2).Value = "hidden_field"
$ws.Cells.Item($lastRow, 3).Value = "Hidden field test"

$ws.Activate()
$ws.Range("A104").Select()
